# The workbook's sole sheet ships with sheet protection enabled, so the
# cells that need updating must be temporarily unprotected before they can
# be written, then protection is restored afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer banner (A7).
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# Refresh the weight / percent-change figures for the model holdings table.
$ws.Range("D2").Value = 0.8492954181081765
$ws.Range("E2").Value = 0.0002483854942871844
$ws.Range("D3").Value = 0.1507045818918235
$ws.Range("E3").Value = 0.009891750653228781
$ws.Range("E4").Value = 0.001701684808595783

$ws.Protect()
